$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert column C (count) for unedited rows 2-214: the engine mis-parses
# leading-space numeric text from the source file as 0, so every count must
# be rewritten explicitly to survive the save round-trip.
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(10, 3).Value = 16
$ws.Cells.Item(11, 3).Value = 19
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(14, 3).Value = 7
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(16, 3).Value = 3
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(19, 3).Value = 32
$ws.Cells.Item(20, 3).Value = 50
$ws.Cells.Item(21, 3).Value = 14
$ws.Cells.Item(22, 3).Value = 6
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(26, 3).Value = 38
$ws.Cells.Item(27, 3).Value = 9
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(29, 3).Value = 3
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(31, 3).Value = 3
$ws.Cells.Item(32, 3).Value = 3
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(36, 3).Value = 3
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(38, 3).Value = 7
$ws.Cells.Item(39, 3).Value = 36
$ws.Cells.Item(40, 3).Value = 16
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(42, 3).Value = 7
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(44, 3).Value = 16
$ws.Cells.Item(45, 3).Value = 15
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(47, 3).Value = 11
$ws.Cells.Item(48, 3).Value = 1
$ws.Cells.Item(49, 3).Value = 12
$ws.Cells.Item(50, 3).Value = 2
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(52, 3).Value = 14
$ws.Cells.Item(53, 3).Value = 2
$ws.Cells.Item(54, 3).Value = 26
$ws.Cells.Item(55, 3).Value = 1
$ws.Cells.Item(56, 3).Value = 4
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(58, 3).Value = 11
$ws.Cells.Item(59, 3).Value = 15
$ws.Cells.Item(60, 3).Value = 3
$ws.Cells.Item(61, 3).Value = 2
$ws.Cells.Item(62, 3).Value = 1
$ws.Cells.Item(63, 3).Value = 1
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(65, 3).Value = 2
$ws.Cells.Item(66, 3).Value = 2
$ws.Cells.Item(67, 3).Value = 13
$ws.Cells.Item(68, 3).Value = 4
$ws.Cells.Item(69, 3).Value = 3
$ws.Cells.Item(70, 3).Value = 4
$ws.Cells.Item(71, 3).Value = 7
$ws.Cells.Item(72, 3).Value = 3
$ws.Cells.Item(73, 3).Value = 19
$ws.Cells.Item(74, 3).Value = 27
$ws.Cells.Item(75, 3).Value = 5
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(77, 3).Value = 2
$ws.Cells.Item(78, 3).Value = 3
$ws.Cells.Item(79, 3).Value = 32
$ws.Cells.Item(80, 3).Value = 1
$ws.Cells.Item(81, 3).Value = 6
$ws.Cells.Item(82, 3).Value = 3
$ws.Cells.Item(83, 3).Value = 1
$ws.Cells.Item(84, 3).Value = 2
$ws.Cells.Item(85, 3).Value = 2
$ws.Cells.Item(86, 3).Value = 3
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(88, 3).Value = 1
$ws.Cells.Item(89, 3).Value = 2
$ws.Cells.Item(90, 3).Value = 11
$ws.Cells.Item(91, 3).Value = 21
$ws.Cells.Item(92, 3).Value = 2
$ws.Cells.Item(93, 3).Value = 10
$ws.Cells.Item(94, 3).Value = 1
$ws.Cells.Item(95, 3).Value = 2
$ws.Cells.Item(96, 3).Value = 7
$ws.Cells.Item(97, 3).Value = 12
$ws.Cells.Item(98, 3).Value = 1
$ws.Cells.Item(99, 3).Value = 3
$ws.Cells.Item(100, 3).Value = 5
$ws.Cells.Item(101, 3).Value = 1
$ws.Cells.Item(102, 3).Value = 2
$ws.Cells.Item(103, 3).Value = 7
$ws.Cells.Item(104, 3).Value = 3
$ws.Cells.Item(105, 3).Value = 3
$ws.Cells.Item(106, 3).Value = 14
$ws.Cells.Item(107, 3).Value = 3
$ws.Cells.Item(108, 3).Value = 1
$ws.Cells.Item(109, 3).Value = 32
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(111, 3).Value = 7
$ws.Cells.Item(112, 3).Value = 4
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(114, 3).Value = 3
$ws.Cells.Item(115, 3).Value = 7
$ws.Cells.Item(116, 3).Value = 5
$ws.Cells.Item(117, 3).Value = 2
$ws.Cells.Item(118, 3).Value = 5
$ws.Cells.Item(119, 3).Value = 1
$ws.Cells.Item(120, 3).Value = 2
$ws.Cells.Item(121, 3).Value = 1
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(123, 3).Value = 3
$ws.Cells.Item(124, 3).Value = 5
$ws.Cells.Item(125, 3).Value = 22
$ws.Cells.Item(126, 3).Value = 1
$ws.Cells.Item(127, 3).Value = 1
$ws.Cells.Item(128, 3).Value = 1
$ws.Cells.Item(129, 3).Value = 9
$ws.Cells.Item(130, 3).Value = 6
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(132, 3).Value = 5
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(134, 3).Value = 2
$ws.Cells.Item(135, 3).Value = 1
$ws.Cells.Item(136, 3).Value = 13
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(138, 3).Value = 60
$ws.Cells.Item(139, 3).Value = 10
$ws.Cells.Item(140, 3).Value = 2
$ws.Cells.Item(141, 3).Value = 3
$ws.Cells.Item(142, 3).Value = 1
$ws.Cells.Item(143, 3).Value = 4
$ws.Cells.Item(144, 3).Value = 55
$ws.Cells.Item(145, 3).Value = 1
$ws.Cells.Item(146, 3).Value = 3
$ws.Cells.Item(147, 3).Value = 1
$ws.Cells.Item(148, 3).Value = 1
$ws.Cells.Item(149, 3).Value = 3
$ws.Cells.Item(150, 3).Value = 13
$ws.Cells.Item(151, 3).Value = 1
$ws.Cells.Item(152, 3).Value = 1
$ws.Cells.Item(153, 3).Value = 1
$ws.Cells.Item(154, 3).Value = 1
$ws.Cells.Item(155, 3).Value = 9
$ws.Cells.Item(156, 3).Value = 3
$ws.Cells.Item(157, 3).Value = 8
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(159, 3).Value = 2
$ws.Cells.Item(160, 3).Value = 1
$ws.Cells.Item(161, 3).Value = 32
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(163, 3).Value = 19
$ws.Cells.Item(164, 3).Value = 11
$ws.Cells.Item(165, 3).Value = 16
$ws.Cells.Item(166, 3).Value = 1
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(168, 3).Value = 8
$ws.Cells.Item(169, 3).Value = 1
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(171, 3).Value = 6
$ws.Cells.Item(172, 3).Value = 1
$ws.Cells.Item(173, 3).Value = 4
$ws.Cells.Item(174, 3).Value = 22
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(177, 3).Value = 4
$ws.Cells.Item(178, 3).Value = 2
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(181, 3).Value = 25
$ws.Cells.Item(182, 3).Value = 2
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(184, 3).Value = 2
$ws.Cells.Item(185, 3).Value = 2
$ws.Cells.Item(186, 3).Value = 73
$ws.Cells.Item(187, 3).Value = 2
$ws.Cells.Item(188, 3).Value = 6
$ws.Cells.Item(189, 3).Value = 7
$ws.Cells.Item(190, 3).Value = 3
$ws.Cells.Item(191, 3).Value = 6
$ws.Cells.Item(192, 3).Value = 17
$ws.Cells.Item(193, 3).Value = 2
$ws.Cells.Item(194, 3).Value = 1
$ws.Cells.Item(195, 3).Value = 4
$ws.Cells.Item(196, 3).Value = 8
$ws.Cells.Item(197, 3).Value = 1
$ws.Cells.Item(198, 3).Value = 10
$ws.Cells.Item(199, 3).Value = 2
$ws.Cells.Item(200, 3).Value = 1
$ws.Cells.Item(201, 3).Value = 12
$ws.Cells.Item(202, 3).Value = 10
$ws.Cells.Item(203, 3).Value = 3
$ws.Cells.Item(204, 3).Value = 2
$ws.Cells.Item(205, 3).Value = 1
$ws.Cells.Item(206, 3).Value = 1
$ws.Cells.Item(207, 3).Value = 4
$ws.Cells.Item(208, 3).Value = 2
$ws.Cells.Item(209, 3).Value = 3
$ws.Cells.Item(210, 3).Value = 2
$ws.Cells.Item(211, 3).Value = 10
$ws.Cells.Item(212, 3).Value = 1
$ws.Cells.Item(213, 3).Value = 34
$ws.Cells.Item(214, 3).Value = 8

# Reorder / update existing 2014 rows 215-266 (B, C only; A already "2014")
$ws.Cells.Item(215, 2).Value = 'Iran'
$ws.Cells.Item(215, 3).Value = 2
$ws.Cells.Item(216, 2).Value = 'Israel'
$ws.Cells.Item(216, 3).Value = 8
$ws.Cells.Item(217, 2).Value = 'Italy'
$ws.Cells.Item(217, 3).Value = 1
$ws.Cells.Item(218, 2).Value = 'Ivory Coast (Cote D`Ivoire)'
$ws.Cells.Item(218, 3).Value = 1
$ws.Cells.Item(219, 2).Value = 'Jordan'
$ws.Cells.Item(219, 3).Value = 3
$ws.Cells.Item(220, 2).Value = 'Kazakhstan'
$ws.Cells.Item(220, 3).Value = 1
$ws.Cells.Item(221, 2).Value = 'Kenya'
$ws.Cells.Item(221, 3).Value = 32
$ws.Cells.Item(222, 2).Value = 'Kuala Lumpur'
$ws.Cells.Item(222, 3).Value = 1
$ws.Cells.Item(223, 2).Value = 'Kyrgyzstan'
$ws.Cells.Item(223, 3).Value = 1
$ws.Cells.Item(224, 2).Value = 'Lebanon'
$ws.Cells.Item(224, 3).Value = 3
$ws.Cells.Item(225, 2).Value = 'Liberia'
$ws.Cells.Item(225, 3).Value = 1
$ws.Cells.Item(226, 2).Value = 'Madagascar'
$ws.Cells.Item(226, 3).Value = 1
$ws.Cells.Item(227, 2).Value = 'Malawi'
$ws.Cells.Item(227, 3).Value = 7
$ws.Cells.Item(228, 2).Value = 'Malaysia'
$ws.Cells.Item(228, 3).Value = 8
$ws.Cells.Item(229, 2).Value = 'Massachusetts'
$ws.Cells.Item(229, 3).Value = 1
$ws.Cells.Item(230, 2).Value = 'Metro Manila'
$ws.Cells.Item(230, 3).Value = 1
$ws.Cells.Item(231, 2).Value = 'Morocco'
$ws.Cells.Item(231, 3).Value = 8
$ws.Cells.Item(232, 2).Value = 'Nairobi'
$ws.Cells.Item(232, 3).Value = 1
$ws.Cells.Item(233, 2).Value = 'Nepal'
$ws.Cells.Item(233, 3).Value = 11
$ws.Cells.Item(234, 2).Value = 'Netherlands'
$ws.Cells.Item(234, 3).Value = 1
$ws.Cells.Item(235, 2).Value = 'New South Wales'
$ws.Cells.Item(235, 3).Value = 2
$ws.Cells.Item(236, 2).Value = 'New Zealand'
$ws.Cells.Item(236, 3).Value = 4
$ws.Cells.Item(237, 2).Value = 'Nigeria'
$ws.Cells.Item(237, 3).Value = 40
$ws.Cells.Item(238, 2).Value = 'Noord Holland'
$ws.Cells.Item(238, 3).Value = 1
$ws.Cells.Item(239, 2).Value = 'Pakistan'
$ws.Cells.Item(239, 3).Value = 17
$ws.Cells.Item(240, 2).Value = 'Palestine'
$ws.Cells.Item(240, 3).Value = 4
$ws.Cells.Item(241, 2).Value = 'Papua New Guinea'
$ws.Cells.Item(241, 3).Value = 1
$ws.Cells.Item(242, 2).Value = 'Philippines'
$ws.Cells.Item(242, 3).Value = 8
$ws.Cells.Item(243, 2).Value = 'Punjab'
$ws.Cells.Item(243, 3).Value = 2
$ws.Cells.Item(244, 2).Value = 'Queensland'
$ws.Cells.Item(244, 3).Value = 1
$ws.Cells.Item(245, 2).Value = 'Romania'
$ws.Cells.Item(245, 3).Value = 1
$ws.Cells.Item(246, 2).Value = 'Rwanda'
$ws.Cells.Item(246, 3).Value = 8
$ws.Cells.Item(247, 2).Value = 'Senegal'
$ws.Cells.Item(247, 3).Value = 1
$ws.Cells.Item(248, 2).Value = 'Sierra Leone'
$ws.Cells.Item(248, 3).Value = 1
$ws.Cells.Item(249, 2).Value = 'Singapore'
$ws.Cells.Item(249, 3).Value = 20
$ws.Cells.Item(250, 2).Value = 'Somalia'
$ws.Cells.Item(250, 3).Value = 2
$ws.Cells.Item(251, 2).Value = 'South Africa'
$ws.Cells.Item(251, 3).Value = 12
$ws.Cells.Item(252, 2).Value = 'South Korea'
$ws.Cells.Item(252, 3).Value = 4
$ws.Cells.Item(253, 2).Value = 'Sri Lanka'
$ws.Cells.Item(253, 3).Value = 1
$ws.Cells.Item(254, 2).Value = 'Sudan'
$ws.Cells.Item(254, 3).Value = 2
$ws.Cells.Item(255, 2).Value = 'Syria'
$ws.Cells.Item(255, 3).Value = 1
$ws.Cells.Item(256, 2).Value = 'Tajikistan'
$ws.Cells.Item(256, 3).Value = 2
$ws.Cells.Item(257, 2).Value = 'Tanzania'
$ws.Cells.Item(257, 3).Value = 8
$ws.Cells.Item(258, 2).Value = 'Thailand'
$ws.Cells.Item(258, 3).Value = 7
$ws.Cells.Item(259, 2).Value = 'Tunisia'
$ws.Cells.Item(259, 3).Value = 1
$ws.Cells.Item(260, 2).Value = 'Turkey'
$ws.Cells.Item(260, 3).Value = 1
$ws.Cells.Item(261, 2).Value = 'Uganda'
$ws.Cells.Item(261, 3).Value = 31
$ws.Cells.Item(262, 2).Value = 'Ukraine'
$ws.Cells.Item(262, 3).Value = 1
$ws.Cells.Item(263, 2).Value = 'United Kingdom'
$ws.Cells.Item(263, 3).Value = 1
$ws.Cells.Item(264, 2).Value = 'United States'
$ws.Cells.Item(264, 3).Value = 19
$ws.Cells.Item(265, 2).Value = 'Uttar Pradesh'
$ws.Cells.Item(265, 3).Value = 2
$ws.Cells.Item(266, 2).Value = 'Vietnam'
$ws.Cells.Item(266, 3).Value = 4

# New trailing 2014 rows 267-269 (copy "2014" text cell to preserve shared-string text type)
$ws.Cells.Item(266, 1).Copy($ws.Cells.Item(267, 1))
$ws.Cells.Item(267, 2).Value = 'Yemen'
$ws.Cells.Item(267, 3).Value = 1
$ws.Cells.Item(266, 1).Copy($ws.Cells.Item(268, 1))
$ws.Cells.Item(268, 2).Value = 'Zambia'
$ws.Cells.Item(268, 3).Value = 1
$ws.Cells.Item(266, 1).Copy($ws.Cells.Item(269, 1))
$ws.Cells.Item(269, 2).Value = 'Zimbabwe'
$ws.Cells.Item(269, 3).Value = 5

# Scratch cell to create "2015" as text (avoids numeric auto-detection), then copy into new rows
$scratch = $ws.Cells.Item(2000, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "2015"
$scratch.Style = $ws.Cells.Item(1, 1).Style

# New 2015 rows 270-324
$scratch.Copy($ws.Cells.Item(270, 1))
$ws.Cells.Item(270, 2).Value = 'Afghanistan'
$ws.Cells.Item(270, 3).Value = 2
$scratch.Copy($ws.Cells.Item(271, 1))
$ws.Cells.Item(271, 2).Value = 'Australia'
$ws.Cells.Item(271, 3).Value = 6
$scratch.Copy($ws.Cells.Item(272, 1))
$ws.Cells.Item(272, 2).Value = 'Austria'
$ws.Cells.Item(272, 3).Value = 1
$scratch.Copy($ws.Cells.Item(273, 1))
$ws.Cells.Item(273, 2).Value = 'Bangladesh'
$ws.Cells.Item(273, 3).Value = 7
$scratch.Copy($ws.Cells.Item(274, 1))
$ws.Cells.Item(274, 2).Value = 'Benin'
$ws.Cells.Item(274, 3).Value = 3
$scratch.Copy($ws.Cells.Item(275, 1))
$ws.Cells.Item(275, 2).Value = 'Botswana'
$ws.Cells.Item(275, 3).Value = 1
$scratch.Copy($ws.Cells.Item(276, 1))
$ws.Cells.Item(276, 2).Value = 'Bulgaria'
$ws.Cells.Item(276, 3).Value = 1
$scratch.Copy($ws.Cells.Item(277, 1))
$ws.Cells.Item(277, 2).Value = 'Burkina Faso'
$ws.Cells.Item(277, 3).Value = 1
$scratch.Copy($ws.Cells.Item(278, 1))
$ws.Cells.Item(278, 2).Value = 'Burundi'
$ws.Cells.Item(278, 3).Value = 2
$scratch.Copy($ws.Cells.Item(279, 1))
$ws.Cells.Item(279, 2).Value = 'Cameroon'
$ws.Cells.Item(279, 3).Value = 10
$scratch.Copy($ws.Cells.Item(280, 1))
$ws.Cells.Item(280, 2).Value = 'Canada'
$ws.Cells.Item(280, 3).Value = 1
$scratch.Copy($ws.Cells.Item(281, 1))
$ws.Cells.Item(281, 2).Value = 'China'
$ws.Cells.Item(281, 3).Value = 1
$scratch.Copy($ws.Cells.Item(282, 1))
$ws.Cells.Item(282, 2).Value = 'Colombia'
$ws.Cells.Item(282, 3).Value = 1
$scratch.Copy($ws.Cells.Item(283, 1))
$ws.Cells.Item(283, 2).Value = 'Congo, Democratic Republic of the (Zaire)'
$ws.Cells.Item(283, 3).Value = 2
$scratch.Copy($ws.Cells.Item(284, 1))
$ws.Cells.Item(284, 2).Value = 'Egypt'
$ws.Cells.Item(284, 3).Value = 1
$scratch.Copy($ws.Cells.Item(285, 1))
$ws.Cells.Item(285, 2).Value = 'Ethiopia'
$ws.Cells.Item(285, 3).Value = 2
$scratch.Copy($ws.Cells.Item(286, 1))
$ws.Cells.Item(286, 2).Value = 'France'
$ws.Cells.Item(286, 3).Value = 3
$scratch.Copy($ws.Cells.Item(287, 1))
$ws.Cells.Item(287, 2).Value = 'Ghana'
$ws.Cells.Item(287, 3).Value = 11
$scratch.Copy($ws.Cells.Item(288, 1))
$ws.Cells.Item(288, 2).Value = 'Guinea'
$ws.Cells.Item(288, 3).Value = 1
$scratch.Copy($ws.Cells.Item(289, 1))
$ws.Cells.Item(289, 2).Value = 'India'
$ws.Cells.Item(289, 3).Value = 33
$scratch.Copy($ws.Cells.Item(290, 1))
$ws.Cells.Item(290, 2).Value = 'Indonesia'
$ws.Cells.Item(290, 3).Value = 9
$scratch.Copy($ws.Cells.Item(291, 1))
$ws.Cells.Item(291, 2).Value = 'Jordan'
$ws.Cells.Item(291, 3).Value = 1
$scratch.Copy($ws.Cells.Item(292, 1))
$ws.Cells.Item(292, 2).Value = 'Kenya'
$ws.Cells.Item(292, 3).Value = 18
$scratch.Copy($ws.Cells.Item(293, 1))
$ws.Cells.Item(293, 2).Value = 'Laos'
$ws.Cells.Item(293, 3).Value = 1
$scratch.Copy($ws.Cells.Item(294, 1))
$ws.Cells.Item(294, 2).Value = 'Lebanon'
$ws.Cells.Item(294, 3).Value = 2
$scratch.Copy($ws.Cells.Item(295, 1))
$ws.Cells.Item(295, 2).Value = 'Madagascar'
$ws.Cells.Item(295, 3).Value = 3
$scratch.Copy($ws.Cells.Item(296, 1))
$ws.Cells.Item(296, 2).Value = 'Malawi'
$ws.Cells.Item(296, 3).Value = 3
$scratch.Copy($ws.Cells.Item(297, 1))
$ws.Cells.Item(297, 2).Value = 'Malaysia'
$ws.Cells.Item(297, 3).Value = 1
$scratch.Copy($ws.Cells.Item(298, 1))
$ws.Cells.Item(298, 2).Value = 'Mongolia'
$ws.Cells.Item(298, 3).Value = 2
$scratch.Copy($ws.Cells.Item(299, 1))
$ws.Cells.Item(299, 2).Value = 'Morocco'
$ws.Cells.Item(299, 3).Value = 1
$scratch.Copy($ws.Cells.Item(300, 1))
$ws.Cells.Item(300, 2).Value = 'Mozambique'
$ws.Cells.Item(300, 3).Value = 1
$scratch.Copy($ws.Cells.Item(301, 1))
$ws.Cells.Item(301, 2).Value = 'Myanmar'
$ws.Cells.Item(301, 3).Value = 1
$scratch.Copy($ws.Cells.Item(302, 1))
$ws.Cells.Item(302, 2).Value = 'Namibia'
$ws.Cells.Item(302, 3).Value = 1
$scratch.Copy($ws.Cells.Item(303, 1))
$ws.Cells.Item(303, 2).Value = 'Nepal'
$ws.Cells.Item(303, 3).Value = 5
$scratch.Copy($ws.Cells.Item(304, 1))
$ws.Cells.Item(304, 2).Value = 'Netherlands'
$ws.Cells.Item(304, 3).Value = 2
$scratch.Copy($ws.Cells.Item(305, 1))
$ws.Cells.Item(305, 2).Value = 'Nigeria'
$ws.Cells.Item(305, 3).Value = 37
$scratch.Copy($ws.Cells.Item(306, 1))
$ws.Cells.Item(306, 2).Value = 'Pakistan'
$ws.Cells.Item(306, 3).Value = 14
$scratch.Copy($ws.Cells.Item(307, 1))
$ws.Cells.Item(307, 2).Value = 'Philippines'
$ws.Cells.Item(307, 3).Value = 7
$scratch.Copy($ws.Cells.Item(308, 1))
$ws.Cells.Item(308, 2).Value = 'Rwanda'
$ws.Cells.Item(308, 3).Value = 5
$scratch.Copy($ws.Cells.Item(309, 1))
$ws.Cells.Item(309, 2).Value = 'Serbia'
$ws.Cells.Item(309, 3).Value = 1
$scratch.Copy($ws.Cells.Item(310, 1))
$ws.Cells.Item(310, 2).Value = 'Singapore'
$ws.Cells.Item(310, 3).Value = 4
$scratch.Copy($ws.Cells.Item(311, 1))
$ws.Cells.Item(311, 2).Value = 'Somalia'
$ws.Cells.Item(311, 3).Value = 3
$scratch.Copy($ws.Cells.Item(312, 1))
$ws.Cells.Item(312, 2).Value = 'South Africa'
$ws.Cells.Item(312, 3).Value = 8
$scratch.Copy($ws.Cells.Item(313, 1))
$ws.Cells.Item(313, 2).Value = 'Sri Lanka'
$ws.Cells.Item(313, 3).Value = 4
$scratch.Copy($ws.Cells.Item(314, 1))
$ws.Cells.Item(314, 2).Value = 'Sweden'
$ws.Cells.Item(314, 3).Value = 1
$scratch.Copy($ws.Cells.Item(315, 1))
$ws.Cells.Item(315, 2).Value = 'Tajikistan'
$ws.Cells.Item(315, 3).Value = 1
$scratch.Copy($ws.Cells.Item(316, 1))
$ws.Cells.Item(316, 2).Value = 'Tanzania'
$ws.Cells.Item(316, 3).Value = 8
$scratch.Copy($ws.Cells.Item(317, 1))
$ws.Cells.Item(317, 2).Value = 'Thailand'
$ws.Cells.Item(317, 3).Value = 1
$scratch.Copy($ws.Cells.Item(318, 1))
$ws.Cells.Item(318, 2).Value = 'Turkey'
$ws.Cells.Item(318, 3).Value = 1
$scratch.Copy($ws.Cells.Item(319, 1))
$ws.Cells.Item(319, 2).Value = 'Uganda'
$ws.Cells.Item(319, 3).Value = 30
$scratch.Copy($ws.Cells.Item(320, 1))
$ws.Cells.Item(320, 2).Value = 'United Kingdom'
$ws.Cells.Item(320, 3).Value = 2
$scratch.Copy($ws.Cells.Item(321, 1))
$ws.Cells.Item(321, 2).Value = 'United States'
$ws.Cells.Item(321, 3).Value = 8
$scratch.Copy($ws.Cells.Item(322, 1))
$ws.Cells.Item(322, 2).Value = 'Vietnam'
$ws.Cells.Item(322, 3).Value = 5
$scratch.Copy($ws.Cells.Item(323, 1))
$ws.Cells.Item(323, 2).Value = 'Zambia'
$ws.Cells.Item(323, 3).Value = 2
$scratch.Copy($ws.Cells.Item(324, 1))
$ws.Cells.Item(324, 2).Value = 'Zimbabwe'
$ws.Cells.Item(324, 3).Value = 9

$scratch.ClearContents()
